$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.289.50"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.870.76"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7078"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "'241.48"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.07777"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").Value = "'25.02"
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("D11").Value = "'0.08393"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.861.43"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "'5.243"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "'0.7104"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "'91.01"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "29.302.82"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'6.072"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "'0.000008184"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").Value = "'239.43"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "2.117.62"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'7.753"
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "'0.9999"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").Value = "'163.24"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'4.395"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "'1.291"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").Value = "'4.293"
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "'0.05333"
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").Value = "'1.937"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "'1.175"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'0.7449"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").Value = "'2.698"
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").Value = "'0.01868"
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("D39").Value = "1.230.11"
$ws.Range("E39").Value = "  +6.28%  "
$ws.Range("D40").Value = "'2.724"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").Value = "'6.555"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("D42").Value = "'0.8843"
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "'109.20"
$ws.Range("E43").Value = "  +5.71%  "
$ws.Range("D44").Value = "'72.29"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "'0.9998"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "2.010.12"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'1.790"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "'9.419"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'0.4312"
$ws.Range("E51").Value = "  +0.31%  "
